$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "insert_compound_names" mistakenly produced a stray "Calculation" column
# (C) with a leftover formula, and the "Frecuencia" counts ended up right
# next to the names instead of after the (now wider) Name column. Fix:
#
#   1. Move the "Frecuencia" header + values from column B to column D.
#   2. Drop the erroneous "Calculation" column (C) contents.
#   3. Merge A:C on every data row into a single "Name" cell.

# 1) Cut column B (header + 6 values) and paste it into column D.
$ws.Range("B1:B7").Cut($ws.Range("D1:D7"))

# 2) Clear out the bogus "Calculation" formulas that lived in column C.
$ws.Range("C1:C7").ClearContents()

# 3) Apply a (default) alignment across the whole A:C block in a single
#    call so every row shares one style, then merge each row's A:C cells.
$ws.Range("A1:C7").HorizontalAlignment = 1

$ws.Range("A1:C1").Merge() | Out-Null
$ws.Range("A2:C2").Merge() | Out-Null
$ws.Range("A3:C3").Merge() | Out-Null
$ws.Range("A4:C4").Merge() | Out-Null
$ws.Range("A5:C5").Merge() | Out-Null
$ws.Range("A6:C6").Merge() | Out-Null
$ws.Range("A7:C7").Merge() | Out-Null

$ws.Range("E10").Select() | Out-Null
